$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '64.199.74'
$ws.Range("E2").Value = '  +1.40%  '

# Row 3
$ws.Range("D3").Value = '3.091.85'
$ws.Range("E3").Value = '  +1.05%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '559.96'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.99%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.54'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.05%  '

# Row 7
$ws.Range("E7").Value = '  +0.09%  '

# Row 8
$ws.Range("D8").Value = '3.090.54'
$ws.Range("E8").Value = '  +1.19%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.506'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.69%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.154'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.09%  '

# Row 11
$ws.Range("E11").Value = '  -4.67%  '

# Row 12
$ws.Range("E12").Value = '  +3.90%  '

# Row 13
$ws.Range("E13").Value = '  +0.98%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.19'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.02%  '

# Row 15
$ws.Range("D15").Value = '3.585.02'
$ws.Range("E15").Value = '  +0.83%  '

# Row 16
$ws.Range("D16").Value = '64.220.95'
$ws.Range("E16").Value = '  +1.54%  '

# Row 17
$ws.Range("D17").Value = '3.084.10'
$ws.Range("E17").Value = '  +0.75%  '

# Row 18
$ws.Range("E18").Value = '  +1.46%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.77'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.44%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '480.71'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.46%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.98'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.28%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.674'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.13%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.54'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.32%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.90'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +10.45%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.24'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.78%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.05%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.80'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.98%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.03'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.40%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.08'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.76%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.15%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.30'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.92%  '

# Row 32
$ws.Range("E32").Value = '  -0.50%  '

# Row 33
$ws.Range("E33").Value = '  +1.56%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.61'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.41%  '

# Row 35
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.21'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.74%  '

# Row 36
$ws.Range("B36").Value = 'OKB'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.82'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.59%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '457.76'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.23%  '

# Row 38
$ws.Range("E38").Value = '  +16.79%  '

# Row 39
$ws.Range("E39").Value = '  +2.76%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0823'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.06%  '

# Row 41
$ws.Range("D41").Value = '2.983.95'
$ws.Range("E41").Value = '  -2.82%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.26'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.08%  '

# Row 43
$ws.Range("E43").Value = '  -2.51%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '28.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.54%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.262'
$ws.Range("D45").Style = "Normal"

# Row 46
$ws.Range("B46").Value = 'Fetch.AI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.14'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.70%  '

# Row 47
$ws.Range("B47").Value = 'USDe'
$ws.Range("C47").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.999'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.03%  '

# Row 48
$ws.Range("E48").Value = '  +2.21%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '120.10'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.77%  '

# Row 50
$ws.Range("D50").Value = '0.0₃0517'
$ws.Range("E50").Value = '  +1.50%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.08'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.61%  '
